$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: (empty) -> false
# Assign via a literal-text formula then paste-special as values so the
# cell lands as a plain shared-string ("false") instead of being
# auto-coerced into a Boolean cell type by a direct Value assignment.
$ws.Range("B7").Formula = '="false"'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Date: updated timestamp
$ws.Range("B8").Value = "2025-07-14T12:58:17-03:00"

# Case Sensitive: (empty) -> true
$ws.Range("B15").Formula = '="true"'
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)

$excel.CutCopyMode = 0
